$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear cells that were removed entirely in rows 2-6
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update C and E values for rows 7-19 with new forecast numbers
$values = @{
    7  = @{ C = 0.9006569003772169;  E = 0.4141666450523163 }
    8  = @{ C = 1.298949644080372;   E = 0.7446484179501223 }
    9  = @{ C = 1.485127130420993;   E = 0.8988967199517361 }
    10 = @{ C = 1.938044824544427;   E = 1.197301207077017 }
    11 = @{ C = 1.730502563828185;   E = 1.20920901052266 }
    12 = @{ C = 2.211325510218898;   E = 1.513838358900466 }
    13 = @{ C = 1.095903126316466;   E = 1.063472944477306 }
    14 = @{ C = 0.8137456736830195;  E = 1.30966355756772 }
    15 = @{ C = -1.434438137829841;  E = 0.8159375071586261 }
    16 = @{ C = 1.85385197842538;    E = 1.2808239555127 }
    17 = @{ C = -0.6079479926716203; E = 0.8021760422591839 }
    18 = @{ C = -0.06520462171909491;E = 0.7367476213790747 }
    19 = @{ C = 0.5869668956646645;  E = 0.8208952814083625 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row].C
    $ws.Cells.Item($row, 5).Value = $values[$row].E
}
